# Excel import can now auto-detect property types from the setter and
# resolve lookup-table entities from a status code in the Excel cell.
# Add a new "Status" column (D) with a sample code value ("E").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Header for Status"
$ws.Range("D2").Value = "E"

# Match the best-fit auto-sizing used for the other header columns.
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666

$ws.Range("D2").Select() | Out-Null
